# Generate Report for Handoff
# Adds a new row (row 3) to each of the three worksheets (Overview, zh-cn, de-de)
# for the newly handed-off file 77153fb4-6455-429b-90ed-7258bd6535db.md

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview" (sheet1) — columns: File Name | Path And Name | Extension |
#   Publish URL | zh-cn | de-de | Latest HO Xliff Generate Date
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A3").Value = "77153fb4-6455-429b-90ed-7258bd6535db.md"
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/297e0918e19fb3646876f183b4fc5143b5785be8/e2e/77153fb4-6455-429b-90ed-7258bd6535db.md", "", "", "e2e\77153fb4-6455-429b-90ed-7258bd6535db.md")
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("D3").Value = ""
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-16 02:36:41"
$wsOverview.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G3"))

# ---------------------------------------------------------------------------
# Sheet "zh-cn" (sheet2) — columns: Source File Name | File Extension | Status |
#   Source Path | Priority | Content Duplicate | Latest Handoff File |
#   Latest Handoff Datetime | Latest Target File | Latest Handback File |
#   Latest Handback DateTime | Reference Tokens | To be localized |
#   Dependency From | Has metadata | Error Detail
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/297e0918e19fb3646876f183b4fc5143b5785be8/e2e/77153fb4-6455-429b-90ed-7258bd6535db.md", "", "", "77153fb4-6455-429b-90ed-7258bd6535db.md")
$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "e2e"
$wsZhCn.Range("E3").Value = "ht"
$wsZhCn.Range("F3").Value = "'False"
$wsZhCn.Range("G3").Value = "77153fb4-6455-429b-90ed-7258bd6535db.1e9acbb341234b8a6ca30a36db541af51e706d8e.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-16 02:36:36"
$wsZhCn.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("K3").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("M3").Value = "'True"
$wsZhCn.Range("O3").Value = "'False"

$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.Resize($wsZhCn.Range("A1:P3"))

# ---------------------------------------------------------------------------
# Sheet "de-de" (sheet3) — same columns as zh-cn
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/297e0918e19fb3646876f183b4fc5143b5785be8/e2e/77153fb4-6455-429b-90ed-7258bd6535db.md", "", "", "77153fb4-6455-429b-90ed-7258bd6535db.md")
$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "e2e"
$wsDeDe.Range("E3").Value = "ht"
$wsDeDe.Range("F3").Value = "'False"
$wsDeDe.Range("G3").Value = "77153fb4-6455-429b-90ed-7258bd6535db.1e9acbb341234b8a6ca30a36db541af51e706d8e.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-16 02:36:41"
$wsDeDe.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("K3").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("M3").Value = "'True"
$wsDeDe.Range("O3").Value = "'False"

$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.Resize($wsDeDe.Range("A1:P3"))
